# Apply the Saldo_guide.xlsx update:
#  - Bump the "Dt. Referencia" (column G) date from 2024-09-04 (45539) to 2024-09-05 (45540)
#    for every data row (rows 2-274).
#  - Update "Saldo Previsto" (E) and "Vl. Total" (H) values for the rows whose balances changed.
#  - Rename the sheet to reflect the new export timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Shift every reference date in column G from 45539 to 45540 ---
$lastRow = 274
$dateRange = $ws.Range("G2:G$lastRow")
$dateRange.Value = 45540

# --- 2. Update the rows whose Saldo Previsto / Vl. Total values changed ---
$changes = @{
    6   = 21068.15
    8   = 46298.18
    51  = 35525.25
    109 = 10367.41
    110 = 33289.48
    112 = 7999.71
    113 = 1.02
    120 = 3011.09
    138 = 65880.22
    143 = 117734.25
    205 = 62188.9
}

foreach ($row in $changes.Keys) {
    $value = $changes[$row]
    $ws.Range("E$row").Value = $value
    $ws.Range("H$row").Value = $value
}

# --- 3. Rename the sheet to match the new export run ---
$ws.Name = "IClientBalance-20240905-083137-"
